# use fieldlabelbuilder for complex types
# - Capitalize the generated field-label strings for complex type children
# - Re-apply the base font so the DejaVu Sans font carries an explicit
#   charset (matches the regenerated styles.xml from the fieldlabelbuilder)
# - Update the remembered selection on the frozen pane

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Relabel the complex-type child/relation headers ---
$ws.Range("B1").Value = "Child1 Name"
$ws.Range("C1").Value = "Child2 Name"
$ws.Range("D1").Value = "AnotherRelation Length"

# --- Make sure the default font explicitly carries a charset ---
$usedRange = $ws.UsedRange
$usedRange.Font.Name = "DejaVu Sans"
$usedRange.Font.Charset = 1

# --- Restore the active selection on the frozen (bottom-left) pane ---
$ws.Range("D13").Select() | Out-Null
